$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46, shifting the blank separator row and the
# summary rows (sum [min] / sum [h] / sum [working weeks]) down by one.
$ws.Rows("46:46").Insert()

# Fill in the new data row 46 with another working-hours entry for 2014-3-4.
$ws.Range("A46").Value = 2014
$ws.Range("B46").Value = 3
$ws.Range("C46").Value = 4
$ws.Range("D46").Value = 0.84375
$ws.Range("E46").Value = 0.91666666666666663
$ws.Range("F46").Formula = "=(E46-D46)*24*60"
$ws.Range("G46").Formula = "=F46/60"

# Update the selection to match the new cursor position recorded in the file.
$ws.Range("C47").Select()
